# Applies the "Neue Zauber mit 8 Runen hinzugefügt" edit:
#  1. Fixes "Resistenzschwächung" -> "Resistenz-Schwächung"
#  2. Fixes "Abweh" -> "Abwehr" (missing trailing r)
#  3. Fixes "Erhöhtz" -> "Erhöht"
#  4. Adds a new "Schneidendes Eis" spell to the "7 Runen" list
#  5. Adds "Waffenrost", "Ölwelle" and "Beherrschung" to the "7 Runen" list
#  6. Fleshes out the "8 Runen" list with many additional spells

$d = $word.ActiveDocument

# NB: this PowerShell host does not honour default parameter values, so
# $listLevel is always passed explicitly by callers below.
# $listLevel uses Word's 1-based ListLevelNumber: 1 => <w:ilvl w:val="0"/>,
# 2 => <w:ilvl w:val="1"/> (which is what every sub-bullet here already is).
function Insert-ListParagraphAfter($anchorRange, [string]$text, $listLevel) {
    # Inserts a brand-new list paragraph right after $anchorRange (which must
    # already be collapsed to the end of the preceding paragraph's text),
    # copying that paragraph's list/style formatting, then fills in $text and
    # leaves $anchorRange collapsed at the end of the freshly inserted
    # paragraph so callers can chain further insertions.
    $anchorRange.InsertParagraphAfter()
    $anchorRange.Move(1, 1) | Out-Null
    $anchorRange.InsertAfter($text)
    $anchorRange.Collapse(0)
    $para = $anchorRange.Paragraphs.First
    if ($para.Range.ListFormat.ListLevelNumber -ne $listLevel) {
        $para.Range.ListFormat.ListLevelNumber = $listLevel
    }
}

# 1) "Große Resistenzschwächung" -> "Große Resistenz-Schwächung"
$d.Content.Find.Execute(
    "Große Resistenzschwächung (2x L, 4x St): 150 BP, Alle Gegner, verringert magische Abwehr",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Große Resistenz-Schwächung (2x L, 4x St): 150 BP, Alle Gegner, verringert magische Abwehr",
    2) | Out-Null

# 2) "-30% Abweh" -> "-30% Abwehr"
$d.Content.Find.Execute(
    "Große Aggression (2x Sch, 4x St): 105 BP, Alle Gegner, +30% Angriff, -30% Abweh",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Große Aggression (2x Sch, 4x St): 105 BP, Alle Gegner, +30% Angriff, -30% Abwehr",
    2) | Out-Null

# 3) "Erhöhtz Angriff" -> "Erhöht Angriff"
$d.Content.Find.Execute(
    "Aura der Blutstärke (7x B): 130 BP, Alle Verbündeten, Erhöhtz Angriff um 20%",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Aura der Blutstärke (7x B): 130 BP, Alle Verbündeten, Erhöht Angriff um 20%",
    2) | Out-Null

# 4) New "Schneidendes Eis" entry, right after "Aura der Blutstärke" (7x B)
$rng = $d.Content
$rng.Find.Execute("Aura der Blutstärke (7x B): 130 BP, Alle Verbündeten, Erhöht Angriff um 20%") | Out-Null
$rng.Collapse(0)
Insert-ListParagraphAfter $rng "Schneidendes Eis (7x E): 172 BP, Alle Gegner, Basisstärke 100, Eisschaden, 60% Chance auf Blutung" 2

# 5) New "Waffenrost", "Ölwelle", "Beherrschung" entries after "Großes Verbennen" (7x F)
$rng = $d.Content
$rng.Find.Execute("Großes Verbennen (7x F): 125 BP, Alle Gegner, 70% Chance auf Brennen") | Out-Null
$rng.Collapse(0)
Insert-ListParagraphAfter $rng "Waffenrost (7x L): 93 BP, Alle Gegner, -10% Angriff" 2
Insert-ListParagraphAfter $rng "Ölwelle (7x Sch): 93 BP, Alle Gegner, Verursacht Öl" 2
Insert-ListParagraphAfter $rng "Beherrschung (7x St): 112 BP, Ein Gegner, Verursacht Beherrschung (Greift Verbündete an, Zustand endet, sobald derjenige Schaden erleidet)" 2

# 6) Flesh out "8 Runen":
#    - insert "Frostbrand" between "Gefühl der Rache" and the existing "Feuerbombe"
$rng = $d.Content
$rng.Find.Execute("Gefühl der Rache (8x B): 60 BP, Ein Verbündeter, Erzeugt RP pro Runde") | Out-Null
$rng.Collapse(0)
Insert-ListParagraphAfter $rng "Frostbrand (8x E): 48 BP, Ein Gegner, Basisstärke 90, Eisschaden, Verursacht Frostbrand" 2

#    - insert the long tail of new spells after the existing "Feuerbombe" (8x F)
$rng = $d.Content
$rng.Find.Execute("Feuerbombe (8x F): 55 BP, Ein Gegner, Verursacht „Feuerbombe“") | Out-Null
$rng.Collapse(0)
Insert-ListParagraphAfter $rng "Heiliges Schild (8x L): 86 BP, Ein Verbündeter, absorbiert Licht- und Blitzschaden" 2
Insert-ListParagraphAfter $rng "Großer Schattenmantel (8x Sch): 105 BP, Alle Verbündeten, +50% Licht- und Schattenresistenz, Schutz vor Verdammnis, Schutz vor Gift" 2
Insert-ListParagraphAfter $rng "Schlafnebel (8x St): 101 BP, Alle Gegner, 60% Chance auf Schlaf" 2
Insert-ListParagraphAfter $rng "Eisresistenz-Aura (4x B, 4x E): 97 BP, Alle Verbündeten, +50% Eisresistenz" 2
Insert-ListParagraphAfter $rng "Feuerresistenz-Aura (4x B, 4x F): 90 BP, Alle Verbündeten, +50% Feuerresistenz" 2
Insert-ListParagraphAfter $rng "Lichtresistenz-Aura (4x B, 4x L): 71 BP, Alle Verbündeten, +50% Lichtresistenz" 2
Insert-ListParagraphAfter $rng "Schattenresistenz-Aura (4x B, 4x Sch): 71 BP, Alle Verbündeten, +50% Schattenresistenz" 2
Insert-ListParagraphAfter $rng "Große Eile (4x B, 4x St): 262 BP, Alle Verbündeten, Gewährt zusätzlichen Schlag pro Angriff (Doppelschlag)" 2
Insert-ListParagraphAfter $rng "Todesregen (2x B, 2x E, 2x F, 2x Sch): 105 BP, Alle Gegner, 10% Chance auf Tod" 2
Insert-ListParagraphAfter $rng "Giftregen (2x B, 2x E, 2x F, 2x St): 150 BP, Alle Gegner, Stärke 70 * Runenmacht, Wasserschaden, 40% Chance auf Nass, 40% Chance auf Gift" 2
Insert-ListParagraphAfter $rng "Todesmarkierung (2x B, 2x E, 2x Sch, 2x St): 78 BP, Ein Gegner, sämtlicher Schaden, den das Ziel während des Zustands erhält wird am Ende teilweise erneut zugefügt" 2
Insert-ListParagraphAfter $rng "Großer Schattenmantel (2x B, 2x F, 2x Sch, 2x St): 105 BP, Alle Verbündeten, +50% Licht- und Schattenresistenz, Schutz vor Verdammnis, Schutz vor Gift [Anmerkung: scheint entweder ein Bug zu sein, oder dieser Zauber hat wirklich 2 verschiedene Runenkonstellationen]" 2
Insert-ListParagraphAfter $rng "Heißer Nebel (4x E, 4x F): 153 BP, Alle Gegner, Basisstärke 80, Wasserschaden, Verursacht Nass, 20% Chance auf Blind" 2
Insert-ListParagraphAfter $rng "Kristallsturm (4x E, 4x L): 176 BP, Alle Gegner, Basisstärke 100, Eisschaden, 20% Chance auf Lähmung" 2
Insert-ListParagraphAfter $rng "Schwarze Eisrüstung (4x E, 4x Sch): 33 BP, Ein Verbündeter, +40% Eis- und Schattenresistenz" 2
Insert-ListParagraphAfter $rng "Frostklingen-Aura (4x E, 4x St): 93 BP, Alle Verbündeten, bei Angriff +15% Chance auf Gefroren" 2
Insert-ListParagraphAfter $rng "Engel der Verdammnis (2x E, 2x L, 2x Sch, 2x St): 262 BP, Alle Gegner, Basisstärke 50, Schattenschaden, Gegner sterben nach 5 Runden (verursacht Verdammnis?)" 2
Insert-ListParagraphAfter $rng "Verbrennendes Licht (4x F, 4x L): 176 BP, Alle Gegner, Basisstärke 95, Lichtschaden, 25% Chance auf Brennen, 25% Chance auf Blind" 2
Insert-ListParagraphAfter $rng "Großer Feuerfluch (4x F, 4x Sch): 168 BP, Alle Gegner, Basisstärke 70, Feuerschaden, Verursacht Heilblocker, 40% Chance auf Blind" 2
Insert-ListParagraphAfter $rng "Lavaklingen-Aura (4x F, 4x St): 78 BP, Alle Verbündeten, +20% Chance auf Brennen bei Angriff" 2
Insert-ListParagraphAfter $rng "Prisma-Explosion (2x F, 2x L, 2x Sch, 2x St): 101 BP, 1 RP, Ein Gegner, Basisstärke 50, Blitzschaden, Verursacht zufällige, negative Zustände" 2
Insert-ListParagraphAfter $rng "Todeslicht (4x L, 4x Sch): 187 BP, Alle Gegner, Basisstärke 80, Schattenschaden, 20% Chance auf Blind, 10% Chance auf Tod" 2
Insert-ListParagraphAfter $rng "Magierbann-Klingen-Aura (4x L, 4x St): 105 BP, Alle Verbündeten, +25% Chance auf Stille bei Angriff" 2
Insert-ListParagraphAfter $rng "Dunkelklingen-Aura (4x Sch, 4x St): 101 BP, Alle Verbündeten, +25% Chance auf Blind bei Angriff" 2

Write-Output "done"
